$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# The previously-blank placeholder rows for the day that follows "Jour 6"
# (rows 29:31) are removed - the day's entries (rows 27:28) are immediately
# followed by its "Bilan du jour" subtotal, same as the other days.
$ws.Rows("29:31").Delete()

# Fill in the now-freed block (rows 30:34, still merged B30:B34 as one
# "day" block) with this day's journal entries.

# Entry 1
$ws.Range("A30").Value = "absent"
$ws.Range("B30").Value2 = 45999
$ws.Range("C30").Value2 = 25
$ws.Range("D30").Value = "I was supposed to go to Sebeillon but actually not"
$ws.Range("E30").Value = "Finished"
$ws.Range("F30").Value2 = 0.35069444444444442
$ws.Range("F30").NumberFormat = "h:mm"

# Entry 2
$ws.Range("A31").Value = "installation"
$ws.Range("C31").Value2 = 35
$ws.Range("D31").Value = "Installing visual studio 2026"
$ws.Range("E31").Value = "Finished"
$ws.Range("F31").Value2 = 0.375
$ws.Range("F31").NumberFormat = "h:mm"

# Entry 3
$ws.Range("A32").Value = "Coding"
$ws.Range("C32").Value2 = 25
$ws.Range("D32").Value = "implementing the bosses waves for later (no bosses now)"
$ws.Range("E32").Value = "Finished"
$ws.Range("F32").Value2 = 0.3923611111111111
$ws.Range("F32").NumberFormat = "h:mm"

# Entry 4
$ws.Range("A33").Value = "Coding"
$ws.Range("C33").Value2 = 45
$ws.Range("D33").Value = "Implementing the pause menu"
$ws.Range("E33").Value = "In the work"
$ws.Range("F33").Value2 = 0.43402777777777773
$ws.Range("F33").NumberFormat = "h:mm"

# Entry 5
$ws.Range("A34").Value = "Coding"
$ws.Range("C34").Value2 = 15
$ws.Range("D34").Value = "fixing the broken wave generation"
$ws.Range("E34").Value = "Finished"
$ws.Range("F34").Value2 = 0.44444444444444442
$ws.Range("F34").NumberFormat = "h:mm"

# Update the print area to match the new (shorter) sheet extent.
$ws.PageSetup.PrintArea = "A1:G42"

"done"
